$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (filename) - only row 2
$ws.Range("A2").Value = ' SCRIPT/T01P01A/enter07.ssb'

# Column B (line numbers)
$ws.Range("B2").Value = 227
$ws.Range("B3").Value = 230
$ws.Range("B4").Value = 233
$ws.Range("B5").Value = 236
$ws.Range("B6").Value = 239

# Column C (English strings)
$ws.Range("C2").Value = ' The great [CS:N]Dusknoir[CR]\''s from the\nfuture! Isn\''t that a startling piece of news!'
$ws.Range("C3").Value = ' But if you think about it, a\nPokémon hatching from an Egg...[K] Wouldn\''t that\nbe like coming from the future?'
$ws.Range("C4").Value = ' ...Pardon?[K] You didn\''t\nunderstand that?'
$ws.Range("C5").Value = ' That is true. It isn\''t an easy\nidea to grasp.'
$ws.Range("C6").Value = ' I didn\''t understand it when I\nsaid it, either! ♪[K] Eep! ♪'

# Column D (translated Russian strings)
$ws.Range("D2").Value = ' Великий [CS:N]Даскнуар[CR] прибыл из\nбудущего! Вот так поразительная новость!'
$ws.Range("D3").Value = ' Но если так подумать, то Покемон,\nвылупляющийся из Яйца...[K] Разве это не\nпохоже на пришествие из будущего?'
$ws.Range("D4").Value = ' ...Что?[K] Ты ничего не понимаешь?'
$ws.Range("D5").Value = ' Всё верно. Эту мысль непросто\nпостичь.'
$ws.Range("D6").Value = ' Хотя, когда я это сказала, я\nтоже её не поняла! ♪[K] Хии! ♪'

# Column E (converted/ciphered strings)
$ws.Range("E2").Value = ' Âåìéëéê [CS:N]Äàòëîôàñ[CR] ðñéáúì éè\náôäôþåãï! Âïó óàë ðïñàèéóåìûîàÿ îïâïòóû!'
$ws.Range("E3").Value = ' Îï åòìé óàë ðïäôíàóû, óï Ðïëåíïî,\nâúìôðìÿýþéêòÿ éè Ÿêøà...[K] Ñàèâå üóï îå\nðïöïçå îà ðñéšåòóâéå éè áôäôþåãï?'
$ws.Range("E4").Value = ' ...Œóï?[K] Óú îéœåãï îå ðïîéíàåšû?'
$ws.Range("E5").Value = ' Âòæ âåñîï. Üóô íúòìû îåðñïòóï\nðïòóéœû.'
$ws.Range("E6").Value = ' Öïóÿ, ëïãäà ÿ üóï òëàèàìà, ÿ\nóïçå åæ îå ðïîÿìà! ♪[K] Öéé! ♪'

# Row heights (match authored worksheet)
$ws.Rows.Item(2).RowHeight = 57.6
$ws.Rows.Item(3).RowHeight = 42
$ws.Rows.Item(4).RowHeight = 21.6
$ws.Rows.Item(5).RowHeight = 21.6
$ws.Rows.Item(6).RowHeight = 21.6

# Selection as recorded in the saved workbook
$ws.Range("B1").Select() | Out-Null
